$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.196431
$ws.Range("N2").Value = 0.589293
$ws.Range("O2").Value = 0.09717285149889213
$ws.Range("P2").Value = 0.09717285149889213
$ws.Range("Q2").Value = 0.163631933775
$ws.Range("R2").Value = 1.472687403975
$ws.Range("S2").Value = 0.09717285149889213
$ws.Range("T2").Value = 0.09717285149889213

# Row 3 updates
$ws.Range("M3").Value = 0.4307096666666667
$ws.Range("O3").Value = 0.2130686423127578
$ws.Range("P3").Value = 0.2130686423127578
$ws.Range("Q3").Value = 0.3587919200750001
$ws.Range("R3").Value = 3.229127280675001
$ws.Range("S3").Value = 0.2130686423127578
$ws.Range("T3").Value = 0.2130686423127578

# Row 4 updates
$ws.Range("O4").Value = 0.68975850618835
$ws.Range("P4").Value = 0.68975850618835
$ws.Range("S4").Value = 0.68975850618835
$ws.Range("T4").Value = 0.68975850618835
